$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell C1: "C" -> "A"
$ws.Range("C1").Value = "A"

# Update row 2 values
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0.009"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "-0.256***"

# Update row 3 label and values
$ws.Range("A3").Value = "A Lag"
$ws.Range("B3").Value = "0.357**"
$ws.Range("C3").Value = "-0.865***"

# Remove rows 4 and 5 (Constant, r2_adj) entirely
$ws.Range("A4:C5").Delete()
